# "notizen und anpassung nach open space"
#
# Appends four new paragraphs at the very end of the document (right
# before the closing sectPr), after the existing paragraph that ends in
# "...Gummersbach/Köln.":
#   1. an empty paragraph
#   2. "Wissenschaftliche Recherche(wenn möglich)"
#   3. "Stakeholdernalyse?"
#   4. "Geocatching für neue historische orte"
#
# All new paragraphs inherit the same spacing (before/after 240) and run
# formatting (color 1F497D / themeColor text2, sz/szCs 21) that the
# preceding paragraph already used.

$d = $word.ActiveDocument

# Collapse a range to the very end of the document content (just before
# the final paragraph mark) and type a paragraph break there. This keeps
# the existing runs of the last paragraph completely untouched and
# yields a clean, truly empty paragraph (matching how Word itself saves
# an empty paragraph - no stray run).
$e = $d.Content.End
$r = $d.Range($e - 1, $e - 1)
$r.Text = [char]13

# "Wissenschaftliche Recherche(wenn möglich)"
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Wissenschaftliche Recherche(wenn möglich)"

# "Stakeholdernalyse?"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Stakeholdernalyse?"

# "Geocatching für neue historische orte"
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Geocatching für neue historische orte"
